$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to text format first so numeric-looking values
# (e.g. "239.64") are stored as text, matching the original inline-string cells.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "36.090.45"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "1.924.27"
$ws.Range("E3").Value = "  -5.12%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "239.64"
$ws.Range("E5").Value = "  -3.72%  "
$ws.Range("D6").Value = "0.602"
$ws.Range("E6").Value = "  -5.56%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "55.13"
$ws.Range("E8").Value = "  -11.60%  "
$ws.Range("D9").Value = "0.360"
$ws.Range("E9").Value = "  -8.77%  "
$ws.Range("D10").Value = "55.15"
$ws.Range("E10").Value = "  -4.89%  "
$ws.Range("D11").Value = "0.0816"
$ws.Range("E11").Value = "  +3.62%  "
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").Value = "2.210.38"
$ws.Range("E13").Value = "  -5.13%  "
$ws.Range("D14").Value = "0.805"
$ws.Range("E14").Value = "  -9.57%  "
$ws.Range("D15").Value = "20.67"
$ws.Range("E15").Value = "  -11.89%  "
$ws.Range("D16").Value = "13.06"
$ws.Range("E16").Value = "  -9.12%  "
$ws.Range("D17").Value = "5.14"
$ws.Range("E17").Value = "  -7.79%  "
$ws.Range("D18").Value = "1.929.13"
$ws.Range("E18").Value = "  -4.83%  "
$ws.Range("D19").Value = "35.994.34"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "68.88"
$ws.Range("E20").Value = "  -4.77%  "
$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("E21").Value = "  -3.13%  "
$ws.Range("D22").Value = "225.51"
$ws.Range("E22").Value = "  -4.43%  "
$ws.Range("D23").Value = "4.90"
$ws.Range("E23").Value = "  -8.62%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  -5.32%  "
$ws.Range("E26").Value = "  -4.12%  "
$ws.Range("D27").Value = "9.20"
$ws.Range("E27").Value = "  -6.65%  "
$ws.Range("D28").Value = "162.39"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").Value = "19.00"
$ws.Range("E29").Value = "  -6.16%  "
$ws.Range("D30").Value = "0.115"
$ws.Range("E30").Value = "  -18.08%  "
$ws.Range("D31").Value = "0.116"
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("E32").Value = "  -6.05%  "
$ws.Range("D33").Value = "4.60"
$ws.Range("E33").Value = "  -8.49%  "
$ws.Range("D34").Value = "0.0616"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").Value = "4.21"
$ws.Range("E35").Value = "  -6.99%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "1.80"
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").Value = "5.88"
$ws.Range("E38").Value = "  -10.64%  "
$ws.Range("D39").Value = "2.12"
$ws.Range("E39").Value = "  -10.70%  "
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  -12.35%  "
$ws.Range("D41").Value = "0.0958"
$ws.Range("E41").Value = "  -5.05%  "
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("E43").Value = "  -9.58%  "
$ws.Range("D44").Value = "0.0205"
$ws.Range("E44").Value = "  -4.86%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.332.60"
$ws.Range("E45").Value = "  -2.19%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "15.34"
$ws.Range("E46").Value = "  -9.12%  "
$ws.Range("D47").Value = "1.01"
$ws.Range("E47").Value = "  -10.87%  "
$ws.Range("D48").Value = "86.56"
$ws.Range("D49").Value = "7.11"
$ws.Range("E49").Value = "  -7.57%  "
$ws.Range("D50").Value = "2.79"
$ws.Range("E50").Value = "  -3.99%  "
$ws.Range("D51").Value = "45.15"
$ws.Range("E51").Value = "  +0.11%  "

# Drop the temporary text-number-format so cells keep their original (default) style.
$ws.Range("B2:E51").ClearFormats()
